$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Status" column (D) entirely - it's empty of data already, just header.
$ws.Range("D1").EntireColumn.Delete()

# Update the URL display text for the Google row (C3) while keeping the
# hyperlink target pointing at the original google URL with client=safari.
$ws.Hyperlinks.Item(2).TextToDisplay = "google.com"

# Update selection to match the new workbook state.
$ws.Range("C9").Select()
